$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Column H ("L_bend_before") values for rows 2-15
$ws.Range("H2").Value  = "straight"
$ws.Range("H3").Value  = "straight"
$ws.Range("H4").Value  = "R = 950, L = 26, n = 2"
$ws.Range("H5").Value  = "R = 950, L = 26, n = 2"
$ws.Range("H6").Value  = "R = 100,L = 11"
$ws.Range("H7").Value  = "straight"
$ws.Range("H8").Value  = "R = 900, L = 25"
$ws.Range("H9").Value  = "R = 800, L = 23"
$ws.Range("H10").Value = "R = 800, L = 14"
$ws.Range("H11").Value = "R = 900, L = 52"
$ws.Range("H12").Value = "R = 600, L = 42"
$ws.Range("H13").Value = "R = 700, L = 50 x 2"
$ws.Range("H14").Value = "R = 240, L = 19, n = 12"
$ws.Range("H15").Value = "R = 325, L = 22, n = 5"

# E15 ("hor div") matches D15 ("vert div") value of "any"
$ws.Range("E15").Value = "any"

# Update the active selection shown when the file was saved
$ws.Range("D13").Select()
